# Update cryptocurrency price (D) and 1h volume change (E) columns
# with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.756.65"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.546.11"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'308.95"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").Value = "'101.54"
$ws.Range("E6").Value = "  +4.97%  "
$ws.Range("D7").Value = "'0.571"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'36.15"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "2.935.31"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "'16.04"
$ws.Range("E15").Value = "  +5.87%  "
$ws.Range("D16").Value = "2.548.93"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "42.766.83"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'6.77"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "'12.40"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "'69.39"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'248.52"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").Value = "'2.91"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'2.07"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'26.55"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").Value = "'40.45"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "'10.13"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").Value = "'157.21"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'5.74"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").Value = "'2.09"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("D38").Value = "'18.17"
$ws.Range("E38").Value = "  -5.88%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "'22.61"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").Value = "'4.20"
$ws.Range("E42").Value = "  +10.21%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "'3.29"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "1.986.58"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "'9.01"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").Value = "2.791.57"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "'81.35"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "'73.47"
$ws.Range("E51").Value = "  -2.03%  "
